# tickets_pideu.xlsx - v.10.9
# "new visual ticket "fixing" and "confirm""
#
# Adds a new "Tiempo de Reparación" column (I) and appends the 2024-05-16
# ticket rows (77-85) that were logged after the previous save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------------
$ws.Range("I1").Value = "Tiempo de Reparación"

# --- Force column A (date strings) to stay TEXT instead of auto-converting -
# to date serials, matching the existing rows 3-76 which store "yyyy-mm-dd"
# as plain text. ClearFormats() afterwards drops the temporary "@" number
# format so the cells end up with the same (default) style as their peers.
$dateRange = $ws.Range("A77:A85")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2024-05-16"
$dateRange.ClearFormats()

# --- Row 77 ------------------------------------------------------------------
$ws.Range("B77").Value = "09:30:59"
$ws.Range("C77").Value = "Fallo tornillo"
$ws.Range("D77").Value = "-"
$ws.Range("E77").Value = "-"
$ws.Range("F77").Value = "-"
$ws.Range("G77").Value = "-"
$ws.Range("H77").Value = "09:31:12"

# --- Row 78 ------------------------------------------------------------------
$ws.Range("B78").Value = "09:31:18"
$ws.Range("C78").Value = "Fallo en elevador"
$ws.Range("D78").Value = "-"
$ws.Range("E78").Value = "-"
$ws.Range("F78").Value = "-"
$ws.Range("G78").Value = "-"
$ws.Range("H78").Value = "09:31:26"

# --- Row 79 ------------------------------------------------------------------
$ws.Range("B79").Value = "09:38:57"
$ws.Range("C79").Value = "-"
$ws.Range("D79").Value = "Tornillo atascado en tolva"
$ws.Range("E79").Value = "-"
$ws.Range("F79").Value = "-"
$ws.Range("G79").Value = "-"
$ws.Range("H79").Value = "09:39:06"
$ws.Range("I79").Value = "0:00:09"

# --- Row 80 ------------------------------------------------------------------
$ws.Range("B80").Value = "09:39:09"
$ws.Range("C80").Value = "-"
$ws.Range("D80").Value = "Fallo etiqueta"
$ws.Range("E80").Value = "-"
$ws.Range("F80").Value = "-"
$ws.Range("G80").Value = "-"
$ws.Range("H80").Value = "09:39:21"
$ws.Range("I80").Value = "0:00:12"

# --- Row 81 ------------------------------------------------------------------
$ws.Range("B81").Value = "09:39:24"
$ws.Range("C81").Value = "-"
$ws.Range("D81").Value = "AOI (fallo etiqueta)"
$ws.Range("E81").Value = "-"
$ws.Range("F81").Value = "-"
$ws.Range("G81").Value = "-"
$ws.Range("H81").Value = "09:39:38"
$ws.Range("I81").Value = "0:00:14"

# --- Row 82 ------------------------------------------------------------------
$ws.Range("B82").Value = "09:49:49"
$ws.Range("C82").Value = "-"
$ws.Range("D82").Value = "No detecta presencia power CP"
$ws.Range("E82").Value = "-"
$ws.Range("F82").Value = "-"
$ws.Range("G82").Value = "-"
$ws.Range("H82").Value = "09:50:01"
$ws.Range("I82").Value = "0:00:12"

# --- Row 83 ------------------------------------------------------------------
$ws.Range("B83").Value = "09:49:57"
$ws.Range("C83").Value = "-"
$ws.Range("D83").Value = "AOI (malla)"
$ws.Range("E83").Value = "-"
$ws.Range("F83").Value = "-"
$ws.Range("G83").Value = "-"
$ws.Range("H83").Value = "09:50:03"
$ws.Range("I83").Value = "0:00:06"

# --- Row 84 ------------------------------------------------------------------
$ws.Range("B84").Value = "09:50:18"
$ws.Range("C84").Value = "-"
$ws.Range("D84").Value = "Cámara no detecta busbar"
$ws.Range("E84").Value = "-"
$ws.Range("F84").Value = "-"
$ws.Range("G84").Value = "-"
$ws.Range("H84").Value = "09:50:27"
$ws.Range("I84").Value = "0:00:09"

# --- Row 85 ------------------------------------------------------------------
$ws.Range("B85").Value = "09:50:39"
$ws.Range("C85").Value = "-"
$ws.Range("D85").Value = "Cámara no detecta Top cover"
$ws.Range("E85").Value = "-"
$ws.Range("F85").Value = "-"
$ws.Range("G85").Value = "-"
$ws.Range("H85").Value = "09:50:45"
$ws.Range("I85").Value = "0:00:06"
